$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A198").Value = "2023-12-11 17:15:45"
$ws.Range("B198").Value = 0.0004

$ws.Range("A199").Value = "2023-12-11 17:16:17"
$ws.Range("B199").Value = 0.0022

$ws.Range("A200").Value = "2023-12-11 17:16:23"
$ws.Range("B200").Value = 0.0004
